$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.283.52'
$ws.Range("E2").Value = '  -0.84%  '
$ws.Range("D3").Value = '3.445.54'
$ws.Range("E3").Value = '  -1.69%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '608.23'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.91%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.23'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.75%  '
$ws.Range("D7").Value = '3.440.70'
$ws.Range("E7").Value = '  -1.71%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.594'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.15%  '
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.192'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.08%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.03'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.30%  '
$ws.Range("E12").Value = '  -3.03%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '44.30'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.36%  '
$ws.Range("E14").Value = '  -2.03%  '
$ws.Range("D15").Value = '4.002.57'
$ws.Range("E15").Value = '  -1.59%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.16'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.35%  '
$ws.Range("D17").Value = '69.364.85'
$ws.Range("E17").Value = '  -0.87%  '
$ws.Range("D18").Value = '3.438.93'
$ws.Range("E18").Value = '  -1.92%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '579.18'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.37%  '
$ws.Range("E20").Value = '  +1.02%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.15'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.19%  '
$ws.Range("E22").Value = '  -2.74%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.90'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.87%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '96.07'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.51%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '15.16'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.79%  '
$ws.Range("E26").Value = '  -2.54%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.43'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.39%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '32.72'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.42%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.63'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.85%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.84'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.06%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.80'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -6.66%  '
$ws.Range("E33").Value = '  -2.69%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.56'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.05%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '580.31'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -17.59%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.53'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.46%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0471'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.15%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0954'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.42%  '
$ws.Range("E39").Value = '  +0.34%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '56.07'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.58%  '
$ws.Range("E41").Value = '  -1.09%  '
$ws.Range("E42").Value = '  -11.52%  '
$ws.Range("D43").Value = '3.242.36'
$ws.Range("E43").Value = '  -2.65%  '
$ws.Range("D44").Value = '0.0₃0687'
$ws.Range("E44").Value = '  -0.24%  '
$ws.Range("E45").Value = '  -5.20%  '
$ws.Range("E46").Value = '  -3.75%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.77'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.98%  '
$ws.Range("E48").Value = '  -6.14%  '
$ws.Range("E49").Value = '  -2.71%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '134.06'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.69%  '
$ws.Range("E51").Value = '  -0.01%  '
